$d = $word.ActiveDocument

# Locate the "Write Up" title paragraph (paragraph 1 of the document) so we
# know where to anchor the new content. The paragraph right after it is the
# first of a handful of empty (Normal-style) placeholder paragraphs.
$titleParaIndex = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq "Write Up") {
        $titleParaIndex = $i
        break
    }
}

# Anchor on the empty placeholder paragraph that immediately follows the
# title, and insert three new paragraphs right before it. Inserting before
# an already-Normal-styled paragraph makes each new paragraph pick up the
# Normal style automatically (matching the plain, style-less paragraphs in
# the target), instead of inheriting the "Title" style of the paragraph
# that precedes it.
$anchor = $d.Paragraphs($titleParaIndex + 1).Range
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()

# Fill in the text for the three newly-created paragraphs.
$bodyPara1 = $titleParaIndex + 1
$bodyPara2 = $titleParaIndex + 2
$headingPara = $titleParaIndex + 3

$d.Paragraphs($bodyPara1).Range.Text = "This week, we will be styling, or formatting the links in the main navigational menu. This will allow our menu to do some transitional fly-ins, along with some more styling maneuvers."
$d.Paragraphs($bodyPara2).Range.Text = "So, if this sounds at all interesting to you, then please join us for our brand-new article entitled:"
$d.Paragraphs($headingPara).Range.Text = "9 Formatting Main Menu Links"

# The final new paragraph is the new article's title, styled as Heading 1.
$d.Paragraphs($headingPara).Range.Style = "Heading1"
